$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 3735
$ws.Range("I3").Value = 3843
$ws.Range("I4").Value = 907
$ws.Range("I5").Value = 357
$ws.Range("I6").Value = 4333
$ws.Range("I7").Value = 13175
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 147
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I6").Value = 55
$ws.Range("I7").Value = 145
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I2").Value = 20
$ws.Range("I7").Value = 66
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("I3").Value = 14
$ws.Range("I7").Value = 44
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I3").Value = 131
$ws.Range("I7").Value = 423
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I6").Value = 69
$ws.Range("I7").Value = 245
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 123
$ws.Range("I6").Value = 168
$ws.Range("I7").Value = 505
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 125
$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 94
$ws.Range("I3").Value = 79
$ws.Range("I6").Value = 90
$ws.Range("I7").Value = 289
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 115
$ws.Range("I5").Value = 41
$ws.Range("I7").Value = 422
$ws.Range("I8").Value = 800
$ws.Range("I11").Value = 206
$ws.Range("I14").Value = 66
$ws.Range("I17").Value = 15
$ws.Range("I18").Value = 89
$ws.Range("I19").Value = 352
$ws.Range("I20").Value = 323
$ws.Range("I22").Value = 37
$ws.Range("I23").Value = 125
$ws.Range("I27").Value = 124
$ws.Range("I29").Value = 851
$ws.Range("I30").Value = 44
$ws.Range("I31").Value = 125
$ws.Range("I33").Value = 596
$ws.Range("I36").Value = 181
$ws.Range("I37").Value = 423
$ws.Range("I44").Value = 94
$ws.Range("I47").Value = 92
$ws.Range("I49").Value = 109
$ws.Range("I52").Value = 291
$ws.Range("I54").Value = 300
$ws.Range("I57").Value = 54
$ws.Range("I60").Value = 65
$ws.Range("I63").Value = 48
$ws.Range("I64").Value = 116
$ws.Range("I65").Value = 289
$ws.Range("I66").Value = 35
$ws.Range("I67").Value = 505
$ws.Range("I75").Value = 48
$ws.Range("I76").Value = 200
$ws.Range("I77").Value = 71
$ws.Range("I78").Value = 189
$ws.Range("I79").Value = 354
$ws.Range("I80").Value = 49
$ws.Range("I83").Value = 267
$ws.Range("I85").Value = 597
$ws.Range("I86").Value = 80
$ws.Range("I88").Value = 118
$ws.Range("I89").Value = 147
$ws.Range("I91").Value = 161
$ws.Range("I94").Value = 126
$ws.Range("I95").Value = 210
$ws.Range("I96").Value = 145
$ws.Range("I98").Value = 88
$ws.Range("I99").Value = 245
$ws.Range("I101").Value = 13175
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 94
$ws.Range("I6").Value = 49
$ws.Range("I7").Value = 267
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 75
$ws.Range("I3").Value = 79
$ws.Range("I7").Value = 210
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 143
$ws.Range("I3").Value = 215
$ws.Range("I5").Value = 22
$ws.Range("I6").Value = 188
$ws.Range("I7").Value = 596
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I3").Value = 11
$ws.Range("I6").Value = 69
$ws.Range("I7").Value = 109
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 69
$ws.Range("I6").Value = 151
$ws.Range("I7").Value = 300
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 254
$ws.Range("I3").Value = 289
$ws.Range("I4").Value = 39
$ws.Range("I6").Value = 236
$ws.Range("I7").Value = 851
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 101
$ws.Range("I6").Value = 95
$ws.Range("I7").Value = 352
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I2").Value = 33
$ws.Range("I7").Value = 94
$ws = $wb.Worksheets.Item("River North")
$ws.Range("I4").Value = 23
$ws.Range("I7").Value = 200
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 156
$ws.Range("I3").Value = 241
$ws.Range("I5").Value = 19
$ws.Range("I7").Value = 597
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I3").Value = 48
$ws.Range("I6").Value = 77
$ws.Range("I7").Value = 189
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 125
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I3").Value = 55
$ws.Range("I7").Value = 161
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 105
$ws.Range("I4").Value = 22
$ws.Range("I7").Value = 354
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I2").Value = 31
$ws.Range("I7").Value = 116
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I3").Value = 99
$ws.Range("I7").Value = 323
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I2").Value = 26
$ws.Range("I7").Value = 89
$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("I3").Value = 7
$ws.Range("I7").Value = 15
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 58
$ws.Range("I3").Value = 55
$ws.Range("I7").Value = 181
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 81
$ws.Range("I3").Value = 97
$ws.Range("I6").Value = 73
$ws.Range("I7").Value = 291
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I2").Value = 23
$ws.Range("I4").Value = 10
$ws.Range("I6").Value = 70
$ws.Range("I7").Value = 126
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 92
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 88
$ws = $wb.Worksheets.Item("North Center")
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 35
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 93
$ws.Range("I7").Value = 206
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I2").Value = 37
$ws.Range("I7").Value = 115
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 30
$ws.Range("I7").Value = 118
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I3").Value = 221
$ws.Range("I6").Value = 261
$ws.Range("I7").Value = 800
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 41
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 34
$ws.Range("I7").Value = 124
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 80
$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I2").Value = 17
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 48
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("I2").Value = 23
$ws.Range("I7").Value = 54
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I2").Value = 16
$ws.Range("I3").Value = 19
$ws.Range("I7").Value = 65
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I2").Value = 12
$ws.Range("I7").Value = 37
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I2").Value = 18
$ws.Range("I7").Value = 71
$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("I3").Value = 10
$ws.Range("I7").Value = 49
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 125
$ws.Range("I7").Value = 422
